$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition listing) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 427
$ws1.Range("F3").Value = 2917
$ws1.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202410/sxfiGFBq1728715876124.jpeg"
$ws1.Range("F5").Value = 57

# --- Sheet "全部类型" (all types listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 427
$ws4.Range("F7").Value = 2917
$ws4.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202410/sxfiGFBq1728715876124.jpeg"
$ws4.Range("F10").Value = 57
